$d = $word.ActiveDocument

$pkgPre = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgPost = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

$rPr = '<w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>'
$listPPr = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' + $rPr + '</w:pPr>'
$plainPPr = '<w:pPr>' + $rPr + '</w:pPr>'
$boldRPr = '<w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>'

# --- Paragraph: "The string inside the println method ..." (index 6) ---
$p6 = '<w:p>' + $listPPr + `
  '<w:r>' + $rPr + '<w:t xml:space="preserve">The string inside the </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r>' + $rPr + '<w:t>println</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r>' + $rPr + '<w:t xml:space="preserve"> method should be closed by quotations on both ends, only one exists, line 3</w:t></w:r>' + `
  '</w:p>'
$d.Paragraphs.Item(6).Range.InsertXML($pkgPre + '<w:body>' + $p6 + '</w:body>' + $pkgPost)

# --- Paragraph: "printLn is incorrect, it should be println, since the first one is not a method, line 4" (index 7) ---
$p7 = '<w:p>' + $listPPr + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r>' + $rPr + '<w:t>printLn</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r>' + $rPr + '<w:t xml:space="preserve"> is incorrect, it should be </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r>' + $rPr + '<w:t>println</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r>' + $rPr + '<w:t>, since the first one is not a method, line 4</w:t></w:r>' + `
  '</w:p>'
$d.Paragraphs.Item(7).Range.InsertXML($pkgPre + '<w:body>' + $p7 + '</w:body>' + $pkgPost)

# --- Paragraph: "System.Out is incorrect, Out should be lowercased, line 5" (index 9) ---
$p9 = '<w:p>' + $listPPr + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r>' + $rPr + '<w:t>System.Out</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r>' + $rPr + '<w:t xml:space="preserve"> is incorrect, Out should be lowercased, line 5</w:t></w:r>' + `
  '</w:p>'
$d.Paragraphs.Item(9).Range.InsertXML($pkgPre + '<w:body>' + $p9 + '</w:body>' + $pkgPost)

# --- Paragraph: "prints is not a method, either use print or println, line 5" (index 10) PLUS the two
#     new trailing paragraphs are all replaced together so no stray empty paragraph is left behind. ---
$p10 = '<w:p>' + $listPPr + `
  '<w:r>' + $rPr + '<w:t xml:space="preserve">prints is not a method, either use print or </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r>' + $rPr + '<w:t>println</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r>' + $rPr + '<w:t>, line 5</w:t></w:r>' + `
  '</w:p>'

$pBlank = '<w:p>' + $plainPPr + '</w:p>'

$pGit = '<w:p>' + $plainPPr + `
  '<w:r>' + $boldRPr + '<w:t>git log:</w:t></w:r>' + `
  '<w:r>' + $rPr + '<w:t xml:space="preserve"> used to check a repository’s committed  history </w:t></w:r>' + `
  '</w:p>'

$p10Para = $d.Paragraphs.Item(10)
$tailRange = $d.Range($p10Para.Range.Start, $d.Content.End)
$tailRange.InsertXML($pkgPre + '<w:body>' + $p10 + $pBlank + $pGit + '</w:body>' + $pkgPost)
